$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B7").Value = "Logica`nLab`nDC1"
$ws.Range("B8").Value = "Logica`nLab`nDC1"
$ws.Range("B7:B8").Interior.Color = 49407
$ws.Range("B9").Select() | Out-Null
